$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Trim leading/trailing whitespace from mayor-name (column E) entries ---
$ws.Range("E2").Value = "张之政"
$ws.Range("E4").Value = "杨军"
$ws.Range("E54").Value = "杨军"
$ws.Range("E11").Value = "唐毅"
$ws.Range("E13").Value = "姜宏"
$ws.Range("E16").Value = "王旺盛"
$ws.Range("E18").Value = "斯琴毕力格"
$ws.Range("E22").Value = "李明伟"
$ws.Range("E28").Value = "吴群刚"
$ws.Range("E36").Value = "杨林兴"
$ws.Range("E43").Value = "邹瑾"
$ws.Range("E47").Value = "张利"
$ws.Range("E48").Value = "杨玉经"
$ws.Range("E49").Value = "杜延安"
$ws.Range("E52").Value = "陈冰冰"
$ws.Range("E55").Value = "操龙灿"
$ws.Range("E56").Value = "戴启远"
$ws.Range("E58").Value = "许继伟"
$ws.Range("E66").Value = " 张海波"
$ws.Range("E75").Value = "田庆盈"
$ws.Range("E80").Value = "李云峰"
$ws.Range("E87").Value = "高键"
$ws.Range("E93").Value = "朱伟"
$ws.Range("E94").Value = "温国辉"
$ws.Range("E95").Value = "刘吉男"
$ws.Range("E103").Value = "马正勇"
$ws.Range("E107").Value = "吕玉印"
$ws.Range("E114").Value = "吴炜"
$ws.Range("E122").Value = "谭丕创"
$ws.Range("E124").Value = "牙生·司地克"
$ws.Range("E125").Value = "买买提明·卡德"
$ws.Range("E127").Value = "王晖"
$ws.Range("E129").Value = "陈金虎"
$ws.Range("E130").Value = "庄兆林"
$ws.Range("E131").Value = "张宝娟"
$ws.Range("E132").Value = "杜小刚"
$ws.Range("E133").Value = "朱立凡现任泰州市委副书记"
$ws.Range("E134").Value = "陈之常男"
$ws.Range("E135").Value = "曹路宝"
$ws.Range("E137").Value = "方伟"
$ws.Range("E138").Value = "朱晓明"
$ws.Range("E139").Value = "陈云"
$ws.Range("E140").Value = "谢来发"
$ws.Range("E141").Value = "黄喜忠"
$ws.Range("E143").Value = "许南吉"
$ws.Range("E156").Value = "邓沛然"
$ws.Range("E160").Value = "张维亮"
$ws.Range("E166").Value = "袁家健"
$ws.Range("E169").Value = "王登喜"
$ws.Range("E170").Value = "刘宛康"
$ws.Range("E171").Value = "刘尚进"
$ws.Range("E173").Value = "徐衣显"
$ws.Range("E176").Value = "朱是西"
$ws.Range("E180").Value = "毛宏芳"
$ws.Range("E182").Value = "刘忻"
$ws.Range("E184").Value = "王纲"
$ws.Range("E185").Value = "盛阅春"
$ws.Range("E187").Value = "汤飞帆"
$ws.Range("E206").Value = "曹立军"
$ws.Range("E209").Value = "阳卫国"
$ws.Range("E214").Value = "刘事青"
$ws.Range("E216").Value = "郑建新"
$ws.Range("E218").Value = "王军"
$ws.Range("E219").Value = "戴超"
$ws.Range("E229").Value = "梁伟新"
$ws.Range("E234").Value = "林兴禄"
$ws.Range("E237").Value = "宋晓路"
$ws.Range("E240").Value = "黄伟"
$ws.Range("E247").Value = "汤方栋"
$ws.Range("E248").Value = "许桂清"
$ws.Range("E251").Value = "隋显利"
$ws.Range("E252").Value = "于学利"
$ws.Range("E256").Value = "郑光照"
$ws.Range("E271").Value = "何忠华"
$ws.Range("E273").Value = "张子林"
$ws.Range("E276").Value = "李世峰"

# --- Fill in missing mayor_sex (F) / mayor_race (G) values ---
$ws.Range("G121").Value = "汉族"
$ws.Range("F134").Value = "男"
$ws.Range("G138").Value = "汉族"
$ws.Range("F139").Value = "男"
$ws.Range("G139").Value = "汉族"
$ws.Range("G140").Value = "汉族"
$ws.Range("F165").Value = "男"
$ws.Range("G215").Value = "汉族"
$ws.Range("G226").Value = "汉族"
$ws.Range("G230").Value = "汉族"
$ws.Range("G233").Value = "汉族"
$ws.Range("G241").Value = "汉族"
$ws.Range("G243").Value = "汉族"
$ws.Range("F251").Value = "男"
$ws.Range("G251").Value = "汉族"
